$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.597.88'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '2.667.46'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D5").Value = '600.25'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").Value = '156.55'
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +2.53%  '
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("D13").Value = '29.39'
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").Value = '0.0000195'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").Value = '3.145.20'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Value = '65.396.13'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '2.662.32'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("E18").Value = '  -1.85%  '
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").Value = '350.51'
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '69.70'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").Value = '0.0000108'
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '9.67'
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  -3.79%  '
$ws.Range("D27").Value = '0.168'
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("D29").Value = '8.06'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").Value = '538.29'
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("E34").Value = '  +2.36%  '
$ws.Range("D35").Value = '5.45'
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("E36").Value = '  -3.03%  '
$ws.Range("D37").Value = '20.40'
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '158.90'
$ws.Range("E39").Value = '  -2.30%  '
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '42.49'
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = '165.42'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("E45").Value = '  -0.61%  '
$ws.Range("E46").Value = '  -5.36%  '
$ws.Range("D47").Value = '23.01'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").Value = '0.647'
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("D50").Value = '0.0997'
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").Value = '19.95'
$ws.Range("E51").Value = '  +0.97%  '
